# Auto-generated edit script: updates crypto price/volume table to match commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.579.74'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '2.285.22'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '96.13'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '266.76'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.59'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.80'
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = '2.628.44'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.11'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.842'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '2.290.41'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '43.605.61'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  +2.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.96'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('E22').Value = '  +6.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.79'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  -9.57%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.12'
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.14'
$ws.Range('E29').Value = '  +1.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.33'
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.81'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0892'
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.35'
$ws.Range('E34').Value = '  -3.81%  '
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0354'
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('E38').Value = '  -3.23%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.235'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('E43').Value = '  +3.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.24'
$ws.Range('E44').Value = '  +5.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.76'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('E46').Value = '  -4.04%  '
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.19'
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.82'
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.429'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.508.45'
$ws.Range('E51').Value = '  +0.46%  '
